# Implement website tracking / active window tracking:
#  - On the "TrackedApps" sheet, the first tracked entry ("Chrome") is
#    repurposed into a "Wikipedia" entry (same time-limit / flags), and a
#    brand-new "Chrome" row is appended at the end with its own settings.
#  - The "TrackedApps" sheet becomes the active/selected sheet (instead of
#    "Assets"), with the selection resting on the first empty row below the
#    data (A5).

$wb = $excel.ActiveWorkbook

$tracked = $wb.Worksheets.Item("TrackedApps")
$assets  = $wb.Worksheets.Item("Assets")

# Rename the existing "Chrome" row to "Wikipedia" - keep its Time Limit,
# Notification and Force Close values untouched.
$tracked.Range("A2").Value = "Wikipedia"

# Append a brand-new "Chrome" tracked-app row underneath the existing data
# (same 10-minute time limit as "Notepad", notifications/force-close off).
$tracked.Range("A4").Value = "Chrome"
$tracked.Range("B4").Value = $tracked.Range("B3").Value()
$tracked.Range("C4").Value = $false
$tracked.Range("D4").Value = $false

# Make "TrackedApps" the active sheet/tab, with the selection on the next
# free row, and drop the old selection state on "Assets".
$tracked.Select()
$tracked.Range("A5").Select()
